# Commit: "this is my 3rd commit"
# The only meaningful content change is the login/password value stored in
# Sheet1!A2 (a shared string), which changes from "Rakesh1111" to "Rakesh11".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("A2").Value = "Rakesh11"
